$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.400.76'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.39'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.38'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5097'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3920'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07789'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.79'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.107'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.95'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.240'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.481'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.817.78'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001146'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +6.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.46'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06628'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.085'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.432.69'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.25'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.253'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.10'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.034.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '155.45'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.399'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.26'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1101'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.103'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.653'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.653'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07056'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2211'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02324'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.186'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.777'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6261'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.16'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.173'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.386'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.42'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.725'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5878'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.22'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.977'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.194'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06895'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.11%  '
